$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the explicit slide background override so the slide falls back to
# the master/layout background (matches the regenerated deck).
try {
    $s.FollowMasterBackground = 1
} catch {
}

# Spire (re)generates slide XML with an empty cSld @name attribute.
try {
    $s.Name = ""
} catch {
}

# Remove the pre-existing shapes: the title "Text Box 1" textbox and the
# old 2x7 "Table 2" table. The new deck replaces both with a single new
# table.
while ($s.Shapes.Count -gt 0) {
    $s.Shapes.Item(1).Delete()
}

# Build the replacement table: 13 rows x 5 columns.
$tblShape = $s.Shapes.AddTable(13, 5, 85, 80, 550, 374.4)
$tblShape.Name = "New Table"

$tbl = $tblShape.Table
$tbl.Columns.Item(1).Width = 100
$tbl.Columns.Item(2).Width = 100
$tbl.Columns.Item(3).Width = 150
$tbl.Columns.Item(4).Width = 100
$tbl.Columns.Item(5).Width = 100

for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Shape.TextFrame.TextRange.ParagraphFormat.Alignment = 1
    }
}

# Keep the frame's on-slide geometry pinned to the exact target size/position
# regardless of any row-height auto-fit the table engine applied above.
$tblShape.Left = 85
$tblShape.Top = 80
$tblShape.Width = 550
$tblShape.Height = 374.4
